$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds values that look numeric ("311.76") but the
# source workbook stores them verbatim as text (inlineStr), e.g. trailing
# zeros like "1.00"/"2.30" and thousand-dot-separated prices like
# "45.065.36" must survive exactly. Force every Price cell we touch to
# Text format first so Excel's automatic number inference never mangles
# (or silently float-rounds) the literal string on assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.065.36'
$ws.Range("E2").Value = '  +2.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.359.72'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.76'
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.67'
$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("E7").Value = '  -0.41%  '

$ws.Range("E8").Value = '  -0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -2.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.75'
$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.42'
$ws.Range("E12").Value = '  -1.40%  '

$ws.Range("E13").Value = '  +1.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.973'
$ws.Range("E14").Value = '  -3.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.714.51'
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.16'
$ws.Range("E16").Value = '  -1.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.360.93'
$ws.Range("E17").Value = '  +1.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.005.59'
$ws.Range("E18").Value = '  +2.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.44'
$ws.Range("E19").Value = '  +11.04%  '

$ws.Range("E20").Value = '  -4.71%  '

$ws.Range("E21").Value = '  -0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.83'
$ws.Range("E22").Value = '  -1.74%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '257.97'
$ws.Range("E24").Value = '  -3.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").Value = '  +1.41%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.06'
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.17'
$ws.Range("E28").Value = '  -6.20%  '

$ws.Range("E29").Value = '  +1.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0966'
$ws.Range("E30").Value = '  +9.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.27'
$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.20'
$ws.Range("E32").Value = '  -4.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '167.71'
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.97'
$ws.Range("E34").Value = '  +5.00%  '

$ws.Range("E35").Value = '  -1.61%  '

$ws.Range("E36").Value = '  +1.04%  '

$ws.Range("E37").Value = '  -1.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.94'
$ws.Range("E38").Value = '  +4.23%  '

$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("E40").Value = '  -3.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.74'
$ws.Range("E41").Value = '  +1.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.96'
$ws.Range("E42").Value = '  -4.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.889.58'
$ws.Range("E43").Value = '  +13.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.19'
$ws.Range("E44").Value = '  -3.26%  '

$ws.Range("E45").Value = '  -4.13%  '

$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.82'
$ws.Range("E46").Value = '  -4.68%  '

$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '81.22'
$ws.Range("E48").Value = '  +5.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.63'
$ws.Range("E49").Value = '  +8.08%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.15'
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.06'
$ws.Range("E51").Value = '  -3.36%  '
